$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = "Household real disposable income"
$ws.Range("A28").Select() | Out-Null
